$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "SubResponse - Text" action column (F) --------------------
# Column F did not exist before; it is a straight copy of column E's
# formatting (fill/border/alignment/width-ish) for the whole used range,
# after which the cell contents that differ are overwritten.
$ws.Range("E1:E14").Copy($ws.Range("F1:F14"))

# Give the title row (row 6) a uniform style across the newly used columns,
# matching the rest of the row (A6:D6) - formats only, keep existing values.
$ws.Range("A6").Copy()
$ws.Range("E6:F6").PasteSpecial(-4122)

# --- Update header / action text -------------------------------------------
# Row 9 (ACTION bodies) and row 10 (ACTION column headers) for the
# "SubResponse" action got split across two columns (E = set subValue,
# F = set Text + push the whole subResponse onto response).
$ws.Range("E9").Value = 'SubResponse subResponse = new SubResponse();
		subResponse.setSubValue($1);'

$ws.Range("F9").Value = 'subResponse.setText($1);
		logger.info("response before setting subResponse: {}", response);
		logger.info("subResponse: {}", subResponse);
		response.setSubResponse(subResponse);		
		logger.info("response after setting subResponse: {}", response);'

$ws.Range("E10").Value = "SubResponse - subValue"
$ws.Range("F10").Value = "SubResponse - Text"

# --- New data column values (rows 11-14) ------------------------------------
$ws.Range("F11").Value = '"sub text 1"'
$ws.Range("F12").Value = '"sub text 2"'
$ws.Range("F13").Value = '"sub text 3"'
$ws.Range("F14").Value = '"sub text 4"'

# --- Layout adjustments ------------------------------------------------------
# Merge the RuleSet-import header across the new column too.
$ws.Range("C2:E2").UnMerge()
$ws.Range("C2:F2").Merge()

# Column widths: E matches D's width, F gets its own (was previously E's).
$ws.Columns("E").ColumnWidth = 49.42857142857143
$ws.Columns("F").ColumnWidth = 61.142857142857146

# Row 9 grew taller to fit the extra action text.
$ws.Rows(9).RowHeight = 75

# Matches the selection recorded in the saved workbook.
$ws.Range("F9").Select()
